$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2
Set-TextCell $ws "D2" "22.409.39"

# Row 3
Set-TextCell $ws "D3" "1.572.71"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
Set-TextCell $ws "D4" "1.001"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
Set-TextCell $ws "D5" "1.001"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
Set-TextCell $ws "D6" "291.62"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
Set-TextCell $ws "D7" "0.3759"
$ws.Range("E7").Value = "  +2.18%  "

# Row 8
Set-TextCell $ws "D8" "49.87"
$ws.Range("E8").Value = "  +1.51%  "

# Row 9
$ws.Range("E9").Value = "  +1.51%  "

# Row 10
Set-TextCell $ws "D10" "0.07646"
$ws.Range("E10").Value = "  +0.78%  "

# Row 11
Set-TextCell $ws "D11" "1.152"
$ws.Range("E11").Value = "  -1.15%  "

# Row 12
Set-TextCell $ws "D12" "1.001"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
Set-TextCell $ws "D13" "21.21"
$ws.Range("E13").Value = "  +0.50%  "

# Row 14
Set-TextCell $ws "D14" "6.017"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15
Set-TextCell $ws "D15" "6.967"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D16" "1.575.88"
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D17" "0.00001135"
$ws.Range("E17").Value = "  +0.40%  "

# Row 18
Set-TextCell $ws "D18" "90.10"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19
Set-TextCell $ws "D19" "0.06732"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
Set-TextCell $ws "D20" "1.001"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
Set-TextCell $ws "D21" "16.78"
$ws.Range("E21").Value = "  +1.70%  "

# Row 22
Set-TextCell $ws "D22" "6.223"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
Set-TextCell $ws "D23" "12.02"
$ws.Range("E23").Value = "  +0.79%  "

# Row 24
Set-TextCell $ws "D24" "22.411.67"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
Set-TextCell $ws "D25" "2.401"
$ws.Range("E25").Value = "  +0.72%  "

# Row 26
Set-TextCell $ws "D26" "2.674"
$ws.Range("E26").Value = "  -10.39%  "

# Row 27
Set-TextCell $ws "D27" "20.21"
$ws.Range("E27").Value = "  +1.93%  "

# Row 28
Set-TextCell $ws "D28" "146.98"
$ws.Range("E28").Value = "  +1.38%  "

# Row 29
Set-TextCell $ws "D29" "5.013"
$ws.Range("E29").Value = "  +1.00%  "

# Row 30
Set-TextCell $ws "D30" "126.71"
$ws.Range("E30").Value = "  +1.40%  "

# Row 31
Set-TextCell $ws "D31" "1.745.82"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
Set-TextCell $ws "D32" "6.171"
$ws.Range("E32").Value = "  -1.08%  "

# Row 33
Set-TextCell $ws "D33" "1.999"
$ws.Range("E33").Value = "  +0.49%  "

# Row 34
Set-TextCell $ws "D34" "0.9845"
$ws.Range("E34").Value = "  -5.53%  "

# Row 35
Set-TextCell $ws "D35" "10.14"
$ws.Range("E35").Value = "  -1.48%  "

# Row 36
Set-TextCell $ws "D36" "0.08515"
$ws.Range("E36").Value = "  +0.69%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D37" "0.02545"
$ws.Range("E37").Value = "  +0.90%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D38" "1.391"
$ws.Range("E38").Value = "  +12.09%  "

# Row 39
Set-TextCell $ws "D39" "0.2316"
$ws.Range("E39").Value = "  -0.34%  "

# Row 40
Set-TextCell $ws "D40" "0.06580"
$ws.Range("E40").Value = "  +1.25%  "

# Row 41
Set-TextCell $ws "D41" "5.447"
$ws.Range("E41").Value = "  -1.27%  "

# Row 42
Set-TextCell $ws "D42" "11.49"
$ws.Range("E42").Value = "  -2.31%  "

# Row 43
Set-TextCell $ws "D43" "0.6412"
$ws.Range("E43").Value = "  +0.96%  "

# Row 44
Set-TextCell $ws "D44" "14.21"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
Set-TextCell $ws "D45" "0.9997"
$ws.Range("E45").Value = "  -0.03%  "

# Row 46
Set-TextCell $ws "D46" "3.802"
$ws.Range("E46").Value = "  +1.41%  "

# Row 47
Set-TextCell $ws "D47" "0.5985"
$ws.Range("E47").Value = "  +0.29%  "

# Row 48
Set-TextCell $ws "D48" "1.297"
$ws.Range("E48").Value = "  +3.04%  "

# Row 49
Set-TextCell $ws "D49" "2.093"
$ws.Range("E49").Value = "  -1.30%  "

# Row 50
Set-TextCell $ws "D50" "124.89"
$ws.Range("E50").Value = "  +1.25%  "

# Row 51
$ws.Range("E51").Value = "  +0.83%  "
